# The commit swaps the contents of ppt/theme/theme1.xml (the "Integral"
# theme used by the slide master / all slides) and ppt/theme/theme2.xml
# (the "Office Theme" used by the notes master): after the edit,
# theme1.xml carries the "Office Theme" color scheme and theme2.xml
# carries the "Integral" color scheme.
#
# The two themes only differ in their <a:clrScheme> (name + the twelve
# scheme colors) - the font scheme and format scheme are byte-identical
# between them - so re-pointing the slide master's theme colors at the
# "Office Theme" palette reproduces the effective visual swap.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# Office Theme palette (this used to live in ppt/theme/theme2.xml) in
# ThemeColorScheme.Item index order: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $cs.Item($i).RGB = $officeColors[$i - 1]
}
